$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the common background value across the whole data range B2:K21
$ws.Range("B2:K21").Value = -18.71649127088704

# Override the cells that differ from the background value
$ws.Range("C2").Value = 2.108486886295224
$ws.Range("I3").Value = 1.242883883813418
$ws.Range("C4").Value = 2.161180047214988
$ws.Range("F4").Value = 3.326189092262672
$ws.Range("H4").Value = 1.480515394731196
$ws.Range("J4").Value = 0.8721520039782286
$ws.Range("C5").Value = 1.503259832237371
$ws.Range("G5").Value = 2.838719921279058
$ws.Range("B7").Value = 2.442264475115138
$ws.Range("E8").Value = 1.815252835935068
$ws.Range("B9").Value = 3.864435463522697
$ws.Range("I10").Value = 1.705997535159666
$ws.Range("K10").Value = 2.215467494674037
$ws.Range("E11").Value = 2.981554156442675
$ws.Range("G11").Value = 2.834575295238892
$ws.Range("K11").Value = 1.971127891903343
$ws.Range("E13").Value = 2.462239828867738
$ws.Range("J13").Value = 1.669359546972029
$ws.Range("K13").Value = 1.737515773158398
$ws.Range("D14").Value = 4.321924913083995
$ws.Range("K14").Value = 1.981291648663512
$ws.Range("J16").Value = 1.905790278147065
$ws.Range("C17").Value = 1.977299380991897
$ws.Range("H17").Value = 2.044115799205702
$ws.Range("I17").Value = 2.102057361397265
$ws.Range("J17").Value = 2.541247023135382
$ws.Range("H18").Value = 1.999163278973651
$ws.Range("I18").Value = 2.07850913127599
$ws.Range("J18").Value = 2.438500466228492
$ws.Range("H19").Value = 1.631442548005602
$ws.Range("I19").Value = 1.789557661067698
$ws.Range("C20").Value = 0.8582888634890971
$ws.Range("F20").Value = 3.317648428803827
$ws.Range("H20").Value = 1.658706226649484
$ws.Range("I20").Value = 1.261052793782498
$ws.Range("K20").Value = 2.053747556126058
$ws.Range("C21").Value = 1.394793432499138
$ws.Range("E21").Value = 1.619022301185725
$ws.Range("G21").Value = 2.514193567662423
$ws.Range("H21").Value = 1.501841370320068
